# Labels geupdate en downloaden template
#
# Updates the "id" counters on rows 4-6 and rewrites the issue-tracker
# entry in row 7 of the "Issue" sheet (ExportIssues.xlsx) with new data:
#   - Gereed switches from "Nee" to "ja" (issue resolved)
#   - Project/Organisatie/Input/Aard ids reset to 1/2/3/4
#   - Actiehouder "Hisham" -> "Dennis"
#   - Kenmerk "Freeze" -> "Crash"
#   - Issues/Antwoord/Opmerking text updated
#   - Manuren 20 -> 4
#   - Dates updated, Status "Open" -> "Afgehandeld"
#   - id counter 1101 -> 1166 (and 1098/1099/1100 -> 1163/1164/1165)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue")

# --- rows 4-6: bump the running "id" counter in column S ---
$ws.Cells.Item(4, 19).Value = 1163
$ws.Cells.Item(5, 19).Value = 1164
$ws.Cells.Item(6, 19).Value = 1165

# --- row 7: new issue record ---
$ws.Cells.Item(7, 1).Value  = "ja"               # A7  Gereed
$ws.Cells.Item(7, 2).Value  = 1                  # B7  Project_Code
$ws.Cells.Item(7, 3).Value  = 2                  # C7  Organisatie_Code
$ws.Cells.Item(7, 4).Value  = 3                  # D7  Input_Bron
$ws.Cells.Item(7, 5).Value  = 4                  # E7  AardId
$ws.Cells.Item(7, 6).Value  = "Pagina"            # F7  Categorie (unchanged)
$ws.Cells.Item(7, 7).Value  = "Dennis"            # G7  Actiehouder
$ws.Cells.Item(7, 8).Value  = "ja"               # H7  Prioriteit (unchanged)
$ws.Cells.Item(7, 9).Value  = "Crash"             # I7  Kenmerk
$ws.Cells.Item(7, 10).Value = "pagina crashed"    # J7  Issues
$ws.Cells.Item(7, 11).Value = "word aan gewerkt"  # K7  Antwoord
$ws.Cells.Item(7, 12).Value = "Geen haast"        # L7  Opmerking
$ws.Cells.Item(7, 13).Value = "Dave"              # M7  Aangever (unchanged)
$ws.Cells.Item(7, 14).Value = 4                  # N7  Manuren
$ws.Cells.Item(7, 15).Value = "4-11-2018 00:00:00" # O7 Datum ingediend
$ws.Cells.Item(7, 16).Value = "3-11-2018 00:00:00" # P7 Datum gepland
$ws.Cells.Item(7, 17).Value = "3-11-2018 00:00:00" # Q7 Datum gereed
$ws.Cells.Item(7, 18).Value = "Afgehandeld"       # R7  Status
$ws.Cells.Item(7, 19).Value = 1166               # S7  id
